$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JCI")

# Row 4 - Inventory
$ws.Range("C4").Value = 1913000000.0
$ws.Range("D4").Value = 1773000000.0
$ws.Range("E4").Value = 1996000000.0
$ws.Range("F4").Value = 2030000000.0
$ws.Range("G4").Value = 1953000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 3417000000.0
$ws.Range("C14").Value = 3210000000.0
$ws.Range("D14").Value = 3120000000.0
$ws.Range("E14").Value = 3057000000.0
$ws.Range("F14").Value = 3119000000.0
$ws.Range("G14").Value = 3336000000.0

# Row 15 - Accrued Expenses
$ws.Range("B15").Value = 817000000.0

# Row 17 - Current Revenue (Deferred)
$ws.Range("B17").Value = 1710000000.0

# Row 18 - Other current liabilities
$ws.Range("B18").Value = 2352000000.0

# Row 22 - Pension and Post-Retirement Liabilities
$ws.Range("B22").Value = 831000000.0

# Row 23 - Long Term Tax Liability (Deferred)
$ws.Range("D23").Value = -477000000.0

# Row 25 - Non-current Liabilities (Other)
$ws.Range("B25").Value = 5529000000.0

# Row 28 - Additional Paid In Capital
$ws.Range("B28").Value = 17034000000.0

# Row 29 - Common Stock (Net)
$ws.Range("B29").Value = 7000000.0

# Row 30 - Retained Earnings
$ws.Range("B30").Value = 2215000000.0

# Row 31 - Treasury Stock
$ws.Range("B31").Value = 1148000000.0

# Row 32 - Common Equity (Total)
$ws.Range("B32").Value = 18757000000.0

# Row 35 - Shares (Common)
$ws.Range("B35").Value = 716715000.0

# Row 37 - Net Debt
$ws.Range("B37").Value = 5884000000.0

# Row 38 - Total Debt
$ws.Range("B38").Value = 7767000000.0
